$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column I
$ws.Range("I1").Value = "Other found locations"

# Row 2 - Irvine et al. authors reformatted; new source tag
$ws.Range("E2").Value = "[Michael%Irvine%NULL%1,    Daniel%Coombs%NULL%1,    Julianne%Skarha%NULL%1,    Brandon%del Pozo%NULL%1,    Josiah%Rich%NULL%1,    Faye%Taxman%NULL%1,    Traci C.%Green%Traci.c.green@gmail.com%1]"
$ws.Range("I2").Value = "_PMC_Springer"

# Row 3 - Truelove et al. authors reformatted; new source tag
$ws.Range("E3").Value = "[Shaun%Truelove%NULL%1,    Orit%Abrahim%NULL%2,    Orit%Abrahim%NULL%0,    Chiara%Altare%NULL%2,    Chiara%Altare%NULL%0,    Stephen A.%Lauer%NULL%2,    Stephen A.%Lauer%NULL%0,    Krya H.%Grantz%NULL%2,    Krya H.%Grantz%NULL%0,    Andrew S.%Azman%NULL%2,    Andrew S.%Azman%NULL%0,    Paul%Spiegel%NULL%2,    Paul%Spiegel%NULL%0,    Parveen%Parmar%NULL%3,    Parveen%Parmar%NULL%0,    Parveen%Parmar%NULL%0,    NULL%NULL%NULL%0,    NULL%NULL%NULL%0,    NULL%NULL%NULL%0]"
$ws.Range("I3").Value = "_PMC"

# Row 4 - Hariri et al. authors reformatted; ID/ID Format reset to not found/N/A; new source tag
$ws.Range("E4").Value = "[ M.%Hariri%null%1,     H.% Rihawi%null%1,     S.% Safadi%null%1,     M. A.% McGlasson%null%1,     W. % Obaid%null%1]"
$ws.Range("F4").Value = "not found"
$ws.Range("G4").Value = "N/A"
$ws.Range("I4").Value = "_MedBiorxiv"

# Row 5 - Bojorquez et al. authors reformatted; ID/ID Format reset to not found/N/A; new source tag
$ws.Range("E5").Value = "[ I.%Bojorquez%null%1,     C.% Infante%null%1,     I.% Vieitez%null%1,     S.% Larrea%null%1,     C. % Santoro%null%1]"
$ws.Range("F5").Value = "not found"
$ws.Range("G5").Value = "N/A"
$ws.Range("I5").Value = "_MedBiorxiv"

# Row 6 - Hintermeier et al. authors reformatted; new source tag
$ws.Range("E6").Value = "[Maren%Hintermeier%NULL%1,    Hande%Gencer%NULL%1,    Katja%Kajikhina%NULL%1,    Sven%Rohleder%NULL%1,    Claudia%Hövener%NULL%1,    Marie%Tallarek%NULL%1,    Jacob%Spallek%NULL%1,    Kayvan%Bozorgmehr%kayvan.bozorgmehr@uni-bielefeld.de%1]"
$ws.Range("I6").Value = "_PMC_elsevier"

# Row 7 - new source tag only (authors unchanged)
$ws.Range("I7").Value = "_PMC"

# Row 8 - Chew et al. authors reformatted; new source tag
$ws.Range("E8").Value = "[MH%Chew%chew.min.hoe@singhealth.com.sg%1,    F.H.%Koh%NULL%1,    JT%Wu%NULL%1,    S.%Ngaserin%NULL%1,    A.%Ng%NULL%1,    BC%Ong%NULL%1,    V.J.%Lee%NULL%1]"
$ws.Range("I8").Value = "_PMC_elsevier"
